$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,16
$row2[0,0] = 3
$row2[0,1] = 1
$row2[0,2] = 463.732605
$row2[0,3] = 1391.197815
$row2[0,4] = 0.3632113435366598
$row2[0,5] = 0.3632113435366598
$row2[0,6] = 3
$row2[0,7] = 1
$row2[0,8] = 15.35884066666667
$row2[0,9] = 46.076522
$row2[0,10] = 0.1012042817263867
$row2[0,11] = 0.1012042817263867
$row2[0,12] = 7122.395192133269
$row2[0,13] = 64101.55672919942
$row2[0,14] = 0.03675854313750353
$row2[0,15] = 0.03675854313750353
$ws.Range("E2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,16
$row3[0,0] = 3
$row3[0,1] = 1
$row3[0,2] = 463.732605
$row3[0,3] = 1391.197815
$row3[0,4] = 0.3632113435366598
$row3[0,5] = 0.3632113435366598
$row3[0,6] = 3
$row3[0,7] = 1
$row3[0,8] = 50.59256466666667
$row3[0,9] = 151.777694
$row3[0,10] = 0.3333704853712116
$row3[0,11] = 0.3333704853712116
$row3[0,12] = 23461.42180650429
$row3[0,13] = 211152.7962585386
$row3[0,14] = 0.1210839418871462
$row3[0,15] = 0.1210839418871462
$ws.Range("E3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,16
$row4[0,0] = 3
$row4[0,1] = 1
$row4[0,2] = 463.732605
$row4[0,3] = 1391.197815
$row4[0,4] = 0.3632113435366598
$row4[0,5] = 0.3632113435366598
$row4[0,6] = 3
$row4[0,7] = 1
$row4[0,8] = 60.37715666666667
$row4[0,9] = 181.13147
$row4[0,10] = 0.397844271305776
$row4[0,11] = 0.397844271305776
$row4[0,12] = 27998.85614352645
$row4[0,13] = 251989.7052917381
$row4[0,14] = 0.1445015522993343
$row4[0,15] = 0.1445015522993343
$ws.Range("E4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,16
$row5[0,0] = 3
$row5[0,1] = 1
$row5[0,2] = 463.732605
$row5[0,3] = 1391.197815
$row5[0,4] = 0.3632113435366598
$row5[0,5] = 0.3632113435366598
$row5[0,6] = 3
$row5[0,7] = 1
$row5[0,8] = 25.43221733333333
$row5[0,9] = 76.296652
$row5[0,10] = 0.1675809615966257
$row5[0,11] = 0.1675809615966258
$row5[0,12] = 11793.74839491282
$row5[0,13] = 106143.7355542154
$row5[0,14] = 0.06086730621267584
$row5[0,15] = 0.06086730621267585
$ws.Range("E5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,16
$row6[0,0] = 3
$row6[0,1] = 1
$row6[0,2] = 12.24662533333333
$row6[0,3] = 36.739876
$row6[0,4] = 0.009591978638444229
$row6[0,5] = 0.009591978638444227
$row6[0,6] = 3
$row6[0,7] = 1
$row6[0,8] = 15.35884066666667
$row6[0,9] = 46.076522
$row6[0,10] = 0.1012042817263867
$row6[0,11] = 0.1012042817263867
$row6[0,12] = 188.0939671990302
$row6[0,13] = 1692.845704791272
$row6[0,14] = 0.0009707493084385924
$row6[0,15] = 0.0009707493084385924
$ws.Range("E6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,16
$row7[0,0] = 3
$row7[0,1] = 1
$row7[0,2] = 12.24662533333333
$row7[0,3] = 36.739876
$row7[0,4] = 0.009591978638444229
$row7[0,5] = 0.009591978638444227
$row7[0,6] = 3
$row7[0,7] = 1
$row7[0,8] = 50.59256466666667
$row7[0,9] = 151.777694
$row7[0,10] = 0.3333704853712116
$row7[0,11] = 0.3333704853712116
$row7[0,12] = 619.5881841251049
$row7[0,13] = 5576.293657125944
$row7[0,14] = 0.003197682574368445
$row7[0,15] = 0.003197682574368445
$ws.Range("E7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,16
$row8[0,0] = 3
$row8[0,1] = 1
$row8[0,2] = 12.24662533333333
$row8[0,3] = 36.739876
$row8[0,4] = 0.009591978638444229
$row8[0,5] = 0.009591978638444227
$row8[0,6] = 3
$row8[0,7] = 1
$row8[0,8] = 60.37715666666667
$row8[0,9] = 181.13147
$row8[0,10] = 0.397844271305776
$row8[0,11] = 0.397844271305776
$row8[0,12] = 739.4164163886356
$row8[0,13] = 6654.74774749772
$row8[0,14] = 0.003816113751792414
$row8[0,15] = 0.003816113751792413
$ws.Range("E8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,16
$row9[0,0] = 3
$row9[0,1] = 1
$row9[0,2] = 12.24662533333333
$row9[0,3] = 36.739876
$row9[0,4] = 0.009591978638444229
$row9[0,5] = 0.009591978638444227
$row9[0,6] = 3
$row9[0,7] = 1
$row9[0,8] = 25.43221733333333
$row9[0,9] = 76.296652
$row9[0,10] = 0.1675809615966257
$row9[0,11] = 0.1675809615966258
$row9[0,12] = 311.4588370772391
$row9[0,13] = 2803.129533695152
$row9[0,14] = 0.001607433003844777
$row9[0,15] = 0.001607433003844777
$ws.Range("E9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,16
$row10[0,0] = 3
$row10[0,1] = 1
$row10[0,2] = 689.7685036666667
$row10[0,3] = 2069.305511
$row10[0,4] = 0.5402504422695089
$row10[0,5] = 0.5402504422695089
$row10[0,6] = 3
$row10[0,7] = 1
$row10[0,8] = 15.35884066666667
$row10[0,9] = 46.076522
$row10[0,10] = 0.1012042817263867
$row10[0,11] = 0.1012042817263867
$row10[0,12] = 10594.04454470142
$row10[0,13] = 95346.40090231274
$row10[0,14] = 0.05467565796224837
$row10[0,15] = 0.05467565796224837
$ws.Range("E10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,16
$row11[0,0] = 3
$row11[0,1] = 1
$row11[0,2] = 689.7685036666667
$row11[0,3] = 2069.305511
$row11[0,4] = 0.5402504422695089
$row11[0,5] = 0.5402504422695089
$row11[0,6] = 3
$row11[0,7] = 1
$row11[0,8] = 50.59256466666667
$row11[0,9] = 151.777694
$row11[0,10] = 0.3333704853712116
$row11[0,11] = 0.3333704853712116
$row11[0,12] = 34897.15762678574
$row11[0,13] = 314074.4186410717
$row11[0,14] = 0.1801035521613979
$row11[0,15] = 0.1801035521613979
$ws.Range("E11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,16
$row12[0,0] = 3
$row12[0,1] = 1
$row12[0,2] = 689.7685036666667
$row12[0,3] = 2069.305511
$row12[0,4] = 0.5402504422695089
$row12[0,5] = 0.5402504422695089
$row12[0,6] = 3
$row12[0,7] = 1
$row12[0,8] = 60.37715666666667
$row12[0,9] = 181.13147
$row12[0,10] = 0.397844271305776
$row12[0,11] = 0.397844271305776
$row12[0,12] = 41646.26100961458
$row12[0,13] = 374816.3490865312
$row12[0,14] = 0.214935543527336
$row12[0,15] = 0.214935543527336
$ws.Range("E12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,16
$row13[0,0] = 3
$row13[0,1] = 1
$row13[0,2] = 689.7685036666667
$row13[0,3] = 2069.305511
$row13[0,4] = 0.5402504422695089
$row13[0,5] = 0.5402504422695089
$row13[0,6] = 3
$row13[0,7] = 1
$row13[0,8] = 25.43221733333333
$row13[0,9] = 76.296652
$row13[0,10] = 0.1675809615966257
$row13[0,11] = 0.1675809615966258
$row13[0,12] = 17542.3424949388
$row13[0,13] = 157881.0824544492
$row13[0,14] = 0.09053568861852665
$row13[0,15] = 0.09053568861852666
$ws.Range("E13:T13").Value = $row13

$row14 = New-Object 'object[,]' 1,16
$row14[0,0] = 3
$row14[0,1] = 1
$row14[0,2] = 111.00921
$row14[0,3] = 333.02763
$row14[0,4] = 0.08694623555538696
$row14[0,5] = 0.08694623555538696
$row14[0,6] = 3
$row14[0,7] = 1
$row14[0,8] = 15.35884066666667
$row14[0,9] = 46.076522
$row14[0,10] = 0.1012042817263867
$row14[0,11] = 0.1012042817263867
$row14[0,12] = 1704.97276892254
$row14[0,13] = 15344.75492030286
$row14[0,14] = 0.008799331318196157
$row14[0,15] = 0.00879933131819616
$ws.Range("E14:T14").Value = $row14

$row15 = New-Object 'object[,]' 1,16
$row15[0,0] = 3
$row15[0,1] = 1
$row15[0,2] = 111.00921
$row15[0,3] = 333.02763
$row15[0,4] = 0.08694623555538696
$row15[0,5] = 0.08694623555538696
$row15[0,6] = 3
$row15[0,7] = 1
$row15[0,8] = 50.59256466666667
$row15[0,9] = 151.777694
$row15[0,10] = 0.3333704853712116
$row15[0,11] = 0.3333704853712116
$row15[0,12] = 5616.24063552058
$row15[0,13] = 50546.16571968522
$row15[0,14] = 0.02898530874829904
$row15[0,15] = 0.02898530874829904
$ws.Range("E15:T15").Value = $row15

$row16 = New-Object 'object[,]' 1,16
$row16[0,0] = 3
$row16[0,1] = 1
$row16[0,2] = 111.00921
$row16[0,3] = 333.02763
$row16[0,4] = 0.08694623555538696
$row16[0,5] = 0.08694623555538696
$row16[0,6] = 3
$row16[0,7] = 1
$row16[0,8] = 60.37715666666667
$row16[0,9] = 181.13147
$row16[0,10] = 0.397844271305776
$row16[0,11] = 0.397844271305776
$row16[0,12] = 6702.420463612901
$row16[0,13] = 60321.7841725161
$row16[0,14] = 0.03459106172731327
$row16[0,15] = 0.03459106172731327
$ws.Range("E16:T16").Value = $row16

$row17 = New-Object 'object[,]' 1,16
$row17[0,0] = 3
$row17[0,1] = 1
$row17[0,2] = 111.00921
$row17[0,3] = 333.02763
$row17[0,4] = 0.08694623555538696
$row17[0,5] = 0.08694623555538696
$row17[0,6] = 3
$row17[0,7] = 1
$row17[0,8] = 25.43221733333333
$row17[0,9] = 76.296652
$row17[0,10] = 0.1675809615966257
$row17[0,11] = 0.1675809615966258
$row17[0,12] = 2823.21035472164
$row17[0,13] = 25408.89319249476
$row17[0,14] = 0.01457053376157848
$row17[0,15] = 0.01457053376157848
$ws.Range("E17:T17").Value = $row17
